$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 488; all rows 488:549 shift down to 489:550.
$ws.Rows("488:488").Insert()

# Populate the newly inserted row 488 with the new record's data.
$ws.Range("A488").Value = 7
$ws.Range("B488").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C488").Value = "Ñuble"
$ws.Range("D488").Value = 45154
$ws.Range("E488").Value = 16
$ws.Range("F488").Value = 100112008
$ws.Range("G488").Value = "Coliflor"
$ws.Range("H488").Value = "Sin especificar"
$ws.Range("I488").Value = "Primera"
$ws.Range("J488").Value = 300
$ws.Range("K488").Value = 1000
$ws.Range("L488").Value = 1000
$ws.Range("M488").Value = 1000
$ws.Range("N488").Value = "$/unidad"
$ws.Range("O488").Value = "Provincia de Diguillín"
$ws.Range("P488").Value = 1000
$ws.Range("Q488").Value = 1
$ws.Range("R488").Value = "Hortaliza"
